$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

$wsZhCn.Range("J2").Value = "3974610d-0856-4f76-94fa-99bda592d0c0.eb58e5290ab270ecee8ed2216865a42d49b4320d.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "3974610d-0856-4f76-94fa-99bda592d0c0.eb58e5290ab270ecee8ed2216865a42d49b4320d.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-08-16 00:57:44"
$wsZhCn.Range("K3").Value = "2016-08-16 00:57:44"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bf76a1676251b023d84398a9f8de8dc6785a67d/e2e/3974610d-0856-4f76-94fa-99bda592d0c0.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "3974610d-0856-4f76-94fa-99bda592d0c0.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bf76a1676251b023d84398a9f8de8dc6785a67d/e2e/3974610d-0856-4f76-94fa-99bda592d0c0.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "3974610d-0856-4f76-94fa-99bda592d0c0.md")

$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("I3").Font.Underline = $true
$wsZhCn.Range("I3").Font.Color = 15570276

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZhCn.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe.Range("J2").Value = "3974610d-0856-4f76-94fa-99bda592d0c0.eb58e5290ab270ecee8ed2216865a42d49b4320d.de-de.xlf"
$wsDeDe.Range("J3").Value = "3974610d-0856-4f76-94fa-99bda592d0c0.eb58e5290ab270ecee8ed2216865a42d49b4320d.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-08-16 00:57:51"
$wsDeDe.Range("K3").Value = "2016-08-16 00:57:51"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bf76a1676251b023d84398a9f8de8dc6785a67d/e2e/3974610d-0856-4f76-94fa-99bda592d0c0.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "3974610d-0856-4f76-94fa-99bda592d0c0.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bf76a1676251b023d84398a9f8de8dc6785a67d/e2e/3974610d-0856-4f76-94fa-99bda592d0c0.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "3974610d-0856-4f76-94fa-99bda592d0c0.md")

$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("I3").Font.Underline = $true
$wsDeDe.Range("I3").Font.Color = 15570276

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDeDe.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1666666666667
